# "Generate Report for Handoff"
#
# b.md finishes a fresh handoff generation: its status flips from
# "Handed back: in sync with en-US" to "Ready for handoff" on the
# Overview sheet as well as on each per-locale sheet, a new handoff
# package (*.xlf) + timestamp is recorded, the "Content Duplicate" flag
# clears, and a version-mismatch warning is written into "Error Detail"
# (whose column is widened to fit the longer text).

$wb = $excel.ActiveWorkbook

$readyStatus = "Ready for handoff"
$overviewHandoffDate = "2016-08-14 16:49:31"
$zhHandoffDate       = "2016-08-14 16:49:23"
$deHandoffDate       = "2016-08-14 16:49:31"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/e2c2d70909c3cc57d44058303b43819cdf83fd43/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/0bf91e9451b7adbb2351107a287110525dff4604/e2e/b.md."

# Excel's ColumnWidth (character units) is offset from the raw OOXML
# <col width="..."> value by 5/6 of a character - subtract that out so the
# saved width lands on exactly 40.
$errorColumnWidth = 40 - (5 / 6)

# ---- Overview sheet: update the b.md row (row 3) ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $readyStatus      # zh-cn status
$overview.Range("F3").Value = $readyStatus      # de-de status
$overview.Range("G3").Value = $overviewHandoffDate  # Latest HO Xliff Generate Date

# ---- zh-cn sheet: update the b.md row (row 3) ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $readyStatus
$zhcn.Range("F2").Copy($zhcn.Range("F3"))   # Content Duplicate: True -> False (copy keeps it text, not a Boolean)
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = $zhHandoffDate
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = $errorColumnWidth

# ---- de-de sheet: update the b.md row (row 3) ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $readyStatus
$dede.Range("F2").Copy($dede.Range("F3"))   # Content Duplicate: True -> False
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = $deHandoffDate
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = $errorColumnWidth
